$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Accuracy (B) and Loss (C) values for the rows whose training
# results are being discarded in preparation for retraining the model.
$ws.Range("B2:C9").ClearContents()
$ws.Range("B11:C11").ClearContents()
